$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric value updates
$ws.Range("D9").Value = 999
$ws.Range("C10").Value = 765788818
$ws.Range("D10").Value = 666621

# Text update
$ws.Range("E10").Value = "Goo"

# Fill color update on D8 (green, RGB(0,176,80) = FF00B050)
$ws.Range("D8").Interior.Color = 5287936

# Move the selection to D22 to match the saved cursor position
$ws.Range("D22").Select()
